$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to stay text so values like "1.00" or
# "43.460.89" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.460.89'
$ws.Range("E2").Value = '  +1.14%  '

$ws.Range("D3").Value = '2.375.67'
$ws.Range("E3").Value = '  +3.21%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = '310.08'
$ws.Range("E5").Value = '  +0.09%  '

$ws.Range("D6").Value = '105.03'
$ws.Range("E6").Value = '  +5.03%  '

$ws.Range("D7").Value = '0.522'
$ws.Range("E7").Value = '  -2.27%  '

$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  +0.35%  '

$ws.Range("D10").Value = '36.32'
$ws.Range("E10").Value = '  +0.69%  '

$ws.Range("D11").Value = '53.36'
$ws.Range("E11").Value = '  +2.49%  '

$ws.Range("D12").Value = '0.0814'
$ws.Range("E12").Value = '  -0.50%  '

$ws.Range("E13").Value = '  -0.49%  '

$ws.Range("D14").Value = '7.03'
$ws.Range("E14").Value = '  -1.73%  '

$ws.Range("D15").Value = '2.744.56'
$ws.Range("E15").Value = '  +3.17%  '

$ws.Range("D16").Value = '15.69'
$ws.Range("E16").Value = '  +5.57%  '

$ws.Range("D17").Value = '2.373.73'
$ws.Range("E17").Value = '  +2.99%  '

$ws.Range("D18").Value = '0.818'
$ws.Range("E18").Value = '  +2.23%  '

$ws.Range("D19").Value = '43.413.77'
$ws.Range("E19").Value = '  +1.11%  '

$ws.Range("D20").Value = '12.02'
$ws.Range("E20").Value = '  -3.74%  '

$ws.Range("D21").Value = '0.0₃0923'
$ws.Range("E21").Value = '  +0.02%  '

$ws.Range("E22").Value = '  +3.43%  '

$ws.Range("E23").Value = '  +0.37%  '

$ws.Range("D24").Value = '242.32'
$ws.Range("E24").Value = '  +1.17%  '

$ws.Range("E25").Value = '  +2.50%  '

$ws.Range("D26").Value = '2.63'
$ws.Range("E26").Value = '  +0.56%  '

$ws.Range("E27").Value = '  +0.28%  '

$ws.Range("D28").Value = '25.92'
$ws.Range("E28").Value = '  +7.47%  '

$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").Value = '37.08'
$ws.Range("E29").Value = '  -3.75%  '

$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").Value = '9.62'
$ws.Range("E30").Value = '  -0.22%  '

$ws.Range("D31").Value = '2.12'
$ws.Range("E31").Value = '  +0.16%  '

$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").Value = '162.13'
$ws.Range("E32").Value = '  -3.62%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '5.29'
$ws.Range("E33").Value = '  -0.88%  '

$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  -0.03%  '

$ws.Range("B35").Value = 'Celestia'
$ws.Range("C35").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D35").Value = '18.34'
$ws.Range("E35").Value = '  +3.76%  '

$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '2.55'
$ws.Range("E36").Value = '  +6.75%  '

$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").Value = '3.13'
$ws.Range("E37").Value = '  +0.15%  '

$ws.Range("D38").Value = '4.74'
$ws.Range("E38").Value = '  +12.22%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.0744'
$ws.Range("E39").Value = '  +0.93%  '

$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = '1.95'
$ws.Range("E40").Value = '  +6.29%  '

$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '0.106'
$ws.Range("E41").Value = '  +1.13%  '

$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").Value = '0.114'
$ws.Range("E42").Value = '  -1.33%  '

$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D43").Value = '2.49'
$ws.Range("E43").Value = '  +8.55%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '20.08'
$ws.Range("E44").Value = '  +4.89%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.005.77'
$ws.Range("E45").Value = '  +1.96%  '

$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '0.0290'
$ws.Range("E46").Value = '  +0.73%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '3.17'
$ws.Range("E47").Value = '  +5.29%  '

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").Value = '10.42'
$ws.Range("E48").Value = '  +6.06%  '

$ws.Range("B49").Value = 'MultiversX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D49").Value = '58.48'
$ws.Range("E49").Value = '  +6.20%  '

$ws.Range("B50").Value = 'HuobiToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D50").Value = '2.96'
$ws.Range("E50").Value = '  -0.49%  '

$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.577.77'
$ws.Range("E51").Value = '  +1.88%  '
